$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new test-plan rows (26-29) ---
$ws.Range("B26").Value = "min char length for item"
$ws.Range("B27").Value = "UTF-8 chars"
$ws.Range("B28").Value = "Unicode Chars"
$ws.Range("B29").Value = "URL encoded string"

# Row 26 grows to accommodate its new content (matches the other wrapped rows)
$ws.Rows.Item(26).RowHeight = 30

# --- Fix typos / wording in existing cells ---
$ws.Range("B12").Value = "Disabilities access (ADA compliant)"
$ws.Range("C12").Value = "color blindness, hearing impaired, visually impaired"
$ws.Range("B19").Value = "try executable code as an entry"
$ws.Range("B23").Value = "multiple users"
$ws.Range("B24").Value = "multiple sessions"
$ws.Range("E2").Value = "Test passes if UI is displayed, no error code presented"
$ws.Range("E3").Value = "Test passes if UI is displayed, no error code presented across all browsers"
$ws.Range("B6").Value = "Access UI across multiple mobile browsers"

# --- Restore the view to the top of the frozen pane ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B6").Select() | Out-Null
